$wb = $excel.ActiveWorkbook

$wsFBS = $wb.Worksheets.Item("FBS")
$wsOther = $wb.Worksheets.Item("Other")

# Update Timestamp column (AK) on the FBS sheet - all data rows share one value
$wsFBS.Range("AK2:AK44").Value = "2025-09-07T16:23:24.672471"

# Update wind_dir_fg values on FBS sheet (column Q)
$wsFBS.Range("Q10").Value = "NW"
$wsFBS.Range("Q16").Value = "NNE"
$wsFBS.Range("Q29").Value = "NE"
$wsFBS.Range("Q37").Value = "SW"

# Update wind_dir_fg values on Other sheet (column S)
$wsOther.Range("S7").Value = "S"
$wsOther.Range("S9").Value = "NNW"
$wsOther.Range("S11").Value = "NNW"
$wsOther.Range("S17").Value = "E"
$wsOther.Range("S26").Value = "E"
$wsOther.Range("S33").Value = "NE"
$wsOther.Range("S37").Value = "E"
$wsOther.Range("S40").Value = "E"
$wsOther.Range("S41").Value = "SSW"
$wsOther.Range("S45").Value = "NW"
$wsOther.Range("S50").Value = "NW"
